# Applies the Tue Sep 26 15:45:47 UTC 2023 cryptos-list refresh to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.158.32'
$ws.Range("E2").Value = '  -0.41%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.584.03'
$ws.Range("E3").Value = '  -0.19%  '

# Row 4
$ws.Range("E4").Value = '  +0.10%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.11'
$ws.Range("E5").Value = '  +0.95%  '

# Row 6
$ws.Range("E6").Value = '  +0.13%  '

# Row 7
$ws.Range("E7").Value = '  +0.13%  '

# Row 8
$ws.Range("E8").Value = '  -0.18%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0605'
$ws.Range("E9").Value = '  -1.02%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.18'
$ws.Range("E10").Value = '  -2.19%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0845'
$ws.Range("E11").Value = '  +0.21%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.807.06'

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.590.11'
$ws.Range("E13").Value = '  +0.06%  '

# Row 14
$ws.Range("E14").Value = '  -1.68%  '

# Row 15
$ws.Range("E15").Value = '  -0.15%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.99'
$ws.Range("E16").Value = '  -1.07%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.177.34'
$ws.Range("E17").Value = '  -0.29%  '

# Row 18
$ws.Range("E18").Value = '  -0.65%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.32'
$ws.Range("E19").Value = '  -0.94%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '212.94'
$ws.Range("E20").Value = '  +0.33%  '

# Row 21
$ws.Range("E21").Value = '  +0.02%  '

# Row 22
$ws.Range("E22").Value = '  -0.76%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.16'
$ws.Range("E23").Value = '  -0.55%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.93'
$ws.Range("E24").Value = '  +0.66%  '

# Row 25
$ws.Range("E25").Value = '  -0.74%  '

# Row 26
$ws.Range("E26").Value = '  -0.02%  '

# Row 27
$ws.Range("E27").Value = '  -1.02%  '

# Row 28
$ws.Range("E28").Value = '  -0.77%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.10'
$ws.Range("E29").Value = '  -1.48%  '

# Row 30
$ws.Range("E30").Value = '  -2.26%  '

# Row 31
$ws.Range("E31").Value = '  +0.24%  '

# Row 32
$ws.Range("E32").Value = '  -1.48%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.339.80'
$ws.Range("E33").Value = '  +3.98%  '

# Row 34
$ws.Range("E34").Value = '  -2.24%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.43'
$ws.Range("E35").Value = '  -0.11%  '

# Row 36
$ws.Range("E36").Value = '  -1.44%  '

# Row 37
$ws.Range("E37").Value = '  -4.26%  '

# Row 38
$ws.Range("E38").Value = '  -0.40%  '

# Row 39
$ws.Range("E39").Value = '  +0.18%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.78'
$ws.Range("E40").Value = '  +2.81%  '

# Row 41
$ws.Range("E41").Value = '  +0.07%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.947'
$ws.Range("E42").Value = '  -16.85%  '

# Row 43
$ws.Range("E43").Value = '  +0.40%  '

# Row 44
$ws.Range("E44").Value = '  -0.26%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.719.58'
$ws.Range("E45").Value = '  +0.01%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '60.93'
$ws.Range("E46").Value = '  -2.83%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '85.85'
$ws.Range("E47").Value = '  -3.35%  '

# Row 48
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.47'
$ws.Range("E48").Value = '  -1.93%  '

# Row 49
$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0983'
$ws.Range("E49").Value = '  -1.82%  '

# Row 50
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0500'
$ws.Range("E50").Value = '  -0.98%  '

# Row 51
$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.998'
$ws.Range("E51").Value = '  -0.11%  '
